$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.230.26'
$ws.Range('E2').Value = '  -1.35%  '
$ws.Range('D3').Value = '2.952.73'
$ws.Range('E3').Value = '  -1.88%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '587.40'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.65'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.40%  '
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('D8').Value = '2.936.99'
$ws.Range('E8').Value = '  -2.28%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.503'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.45%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.82'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +9.51%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.146'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.49%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.451'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.12%  '
$ws.Range('E13').Value = '  -1.70%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.58'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.47%  '
$ws.Range('E15').Value = '  -0.39%  '
$ws.Range('D16').Value = '3.443.57'
$ws.Range('E16').Value = '  -1.82%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.90'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.92%  '
$ws.Range('D18').Value = '61.254.75'
$ws.Range('E18').Value = '  -1.27%  '
$ws.Range('D19').Value = '2.944.67'
$ws.Range('E19').Value = '  -2.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '432.39'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.75%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.88'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.49%  '
$ws.Range('E22').Value = '  -1.59%  '
$ws.Range('E23').Value = '  -0.99%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '80.66'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.28%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.01'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.37%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.20'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.47%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.95'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.74%  '
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.36'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.19%  '
$ws.Range('E30').Value = '  -0.19%  '
$ws.Range('E31').Value = '  +5.60%  '
$ws.Range('E32').Value = '  -2.76%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '26.99'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.20%  '
$ws.Range('E34').Value = '  -3.10%  '
$ws.Range('E35').Value = '  +0.38%  '
$ws.Range('E36').Value = '  -1.30%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.73'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.37%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '50.05'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.97%  '
$ws.Range('B39').Value = 'dogwifhat'
$ws.Range('C39').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.00'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.71%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.127'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.21%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.04'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.75%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.78'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.24%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.293'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +7.22%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '42.07'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.33%  '
$ws.Range('E45').Value = '  -0.85%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '374.55'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -6.11%  '
$ws.Range('D47').Value = '2.666.21'
$ws.Range('E47').Value = '  -1.77%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '133.25'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.14%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '25.66'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +8.41%  '
$ws.Range('E50').Value = '  -0.03%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.107'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.84%  '
